# Append the "success" item-details row (AWB No / Credit Reference / Date)
# generated by the bulk-upload process onto row 2 of the Success sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AWB No looks fully numeric ("77707350195"), so Excel would normally
# auto-convert it to a number when assigned via .Value. Force it to be
# stored as text first, then clear the number-format override so the
# cell keeps using the default (unstyled) cell style, matching the
# other text cells on the sheet.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "77707350195"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = "CR940126"
$ws.Range("C2").Value = "2026-02-04T13:31:26.532791488"
